{"js": "// Append a new dated log entry (\"25 November 2016 (15mins)\") describing the\n// require-statement change, reusing the two trailing empty \"List Bullet\"\n// placeholder paragraphs that were sitting at the end of the document, and\n// move the trailing `_GoBack` bookmark onto the newly written text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two trailing placeholder paragraphs: both are empty, styled\n// \"List Bullet\", and carry a `<w:numPr><w:numId w:val=\"0\"/></w:numPr>`\n// override (i.e. a \"List Bullet\" paragraph with its bullet switched off).\n// They are the last two paragraphs in the body, immediately preceded by a\n// similar placeholder paragraph that additionally carries explicit\n// indentation (`w:ind`), which we leave untouched.\nlet headingPara = items[items.length - 2];\nlet entryPara = items[items.length - 1];\n\n// Sanity guard: both target paragraphs must be empty. If the document\n// shape differs from what is expected, fall back to the last two empty\n// paragraphs found anywhere in the body.\nif (headingPara.text.trim() !== \"\" || entryPara.text.trim() !== \"\") {\n  const empties = items.filter((p) => p.text.trim() === \"\");\n  headingPara = empties[empties.length - 2];\n  entryPara = empties[empties.length - 1];\n}\n\n// Turn the first placeholder into the new bold date/duration heading,\n// keeping its paragraph-level \"no bullet\" override.\nconst headingOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n  <w:pPr>\n    <w:pStyle w:val=\"ListBullet\"/>\n    <w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"0\"/></w:numPr>\n    <w:rPr><w:b/></w:rPr>\n  </w:pPr>\n  <w:r><w:rPr><w:b/></w:rPr><w:t>25 November 2016 (15mins)</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\nheadingPara.insertOoxml(headingOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// Turn the second placeholder into the plain bulleted log entry (no numPr\n// override here, so it regains the normal \"List Bullet\" bullet glyph).\nconst entryOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\"><pkg:xmlData>\n<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>\n</pkg:xmlData></pkg:part>\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>\n<w:p>\n  <w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>\n  <w:r><w:t>Added the require statement for the Database.js in app.js hoping that it will allow the connection to be set initially when the server is started. Previously the connection was done when the index page was loaded. Seemed to be working fine locally.</w:t></w:r>\n</w:p>\n</w:body></w:document>\n</pkg:xmlData></pkg:part>\n</pkg:package>`;\nentryPara.insertOoxml(entryOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\n// The `_GoBack` bookmark used to sit right after \"...solve it.\" at the end\n// of the LESSON paragraph; Word always re-anchors it to the location of the\n// most recent edit, so drop the old one and re-insert it at the end of the\n// freshly-typed log entry.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst newLastPara = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nconst endRange = newLastPara.getRange(Word.RangeLocation.end);\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Append a new dated log entry (\"25 November 2016 (15mins)\") describing the\n# require-statement change, reusing the two trailing empty \"List Bullet\"\n# placeholder paragraphs that were sitting at the end of the document, and\n# move the trailing `_GoBack` bookmark onto the newly written text.\n\n$d = $word.ActiveDocument\n\n# Locate the two trailing placeholder paragraphs to repurpose: both are\n# empty, styled \"List Bullet\", have no active bullet (numbering switched off\n# via a paragraph-level numId=0 override -> ListFormat.ListType is\n# wdListNoNumbering) and, unlike the placeholder right before them, carry no\n# explicit left indent.\n$n = $d.Paragraphs.Count\n$headingPara = $null\n$entryPara = $null\nfor ($i = $n; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim().Length -eq 0 -and `\n        $p.Style.NameLocal -eq \"List Bullet\" -and `\n        $p.Range.ListFormat.ListType -eq 0 -and `\n        $p.LeftIndent -eq 0) {\n        if ($null -eq $entryPara) {\n            $entryPara = $p\n        } elseif ($null -eq $headingPara) {\n            $headingPara = $p\n            break\n        }\n    }\n}\n\n# Fallback to the last two paragraphs if the expected shape was not found.\nif ($null -eq $headingPara -or $null -eq $entryPara) {\n    $headingPara = $d.Paragraphs.Item($n - 1)\n    $entryPara = $d.Paragraphs.Item($n)\n}\n\n# The `_GoBack` bookmark used to sit right after \"...solve it.\" at the end\n# of the LESSON paragraph; Word always re-anchors it to the location of the\n# most recent edit, so drop the old one now and re-create it after the\n# freshly-typed log entry below.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# --- New bold date/duration heading -----------------------------------\n# Reuses the placeholder \"List Bullet\" paragraph, keeping its \"no bullet\"\n# (numId=0) override, and makes the whole paragraph (including its mark)\n# bold.\n$headingPara.Range.InsertAfter(\"25 November 2016 (15mins)\")\n$headingPara.Range.Font.Bold = 1\n\n# --- New plain bulleted log entry --------------------------------------\n# Re-applying the style drops the paragraph-level numId=0 override, so this\n# paragraph regains the normal \"List Bullet\" bullet glyph.\n$entryPara.Style = \"List Bullet\"\n$entryText = \"Added the require statement for the Database.js in app.js hoping that it will allow the connection to be set initially when the server is started. Previously the connection was done when the index page was loaded. Seemed to be working fine locally.\"\n\n$r = $entryPara.Range\n[void]$r.MoveEnd(1, -1)\n# Type the real text plus one throw-away marker character so we have a\n# *non-degenerate* range to bookmark, then delete just that marker\n# character. Deleting bookmarked text collapses the bookmark down to zero\n# width exactly at that spot, which is where Word normally leaves\n# `_GoBack` after the last edit.\n$r.InsertAfter($entryText + \"#\")\n$markerStart = $r.End - 1\n$markerEnd = $r.End\n$d.Bookmarks.Add(\"_GoBack\", $d.Range($markerStart, $markerEnd))\n$d.Range($markerStart, $markerEnd).Delete()\n"}
